$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.044.11'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '1.598.30'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('D6').Value = '302.10'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').Value = '0.3633'
$ws.Range('E8').Value = '  -0.73%  '
$ws.Range('D9').Value = '50.86'
$ws.Range('E9').Value = '  +4.48%  '
$ws.Range('E10').Value = '  -2.22%  '
$ws.Range('D11').Value = '1.001'
$ws.Range('E11').Value = '  -0.36%  '
$ws.Range('D12').Value = '0.08136'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').Value = '22.33'
$ws.Range('E13').Value = '  -2.41%  '
$ws.Range('D14').Value = '6.564'
$ws.Range('E14').Value = '  -1.07%  '
$ws.Range('D15').Value = '7.346'
$ws.Range('E15').Value = '  -2.94%  '
$ws.Range('D16').Value = '0.00001241'
$ws.Range('E16').Value = '  -1.54%  '
$ws.Range('D17').Value = '1.600.16'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').Value = '92.35'
$ws.Range('E18').Value = '  +0.29%  '
$ws.Range('D19').Value = '0.06849'
$ws.Range('E19').Value = '  +0.55%  '
$ws.Range('E20').Value = '  -1.85%  '
$ws.Range('D21').Value = '6.499'
$ws.Range('E21').Value = '  -1.66%  '
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').Value = '12.98'
$ws.Range('E23').Value = '  -1.11%  '
$ws.Range('D24').Value = '23.045.54'
$ws.Range('E24').Value = '  -0.50%  '
$ws.Range('D25').Value = '2.373'
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('D26').Value = '2.800'
$ws.Range('E26').Value = '  -5.51%  '
$ws.Range('D27').Value = '21.09'
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('D28').Value = '149.05'
$ws.Range('E28').Value = '  -1.29%  '
$ws.Range('D29').Value = '5.252'
$ws.Range('E29').Value = '  +0.44%  '
$ws.Range('D30').Value = '134.94'
$ws.Range('E30').Value = '  +2.03%  '
$ws.Range('D31').Value = '2.377'
$ws.Range('E31').Value = '  -3.63%  '
$ws.Range('D32').Value = '6.735'
$ws.Range('E32').Value = '  -5.28%  '
$ws.Range('D33').Value = '1.774.18'
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('D34').Value = '0.9584'
$ws.Range('E34').Value = '  -1.66%  '
$ws.Range('D35').Value = '0.07508'
$ws.Range('E35').Value = '  -3.09%  '
$ws.Range('D36').Value = '0.02709'
$ws.Range('E36').Value = '  -2.64%  '
$ws.Range('B37').Value = 'FraxShare'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D37').Value = '10.16'
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = '6.197'
$ws.Range('E38').Value = '  -1.06%  '
$ws.Range('E39').Value = '  -1.43%  '
$ws.Range('D40').Value = '0.08814'
$ws.Range('E40').Value = '  -0.68%  '
$ws.Range('D41').Value = '1.356'
$ws.Range('E41').Value = '  -2.09%  '
$ws.Range('D42').Value = '0.7032'
$ws.Range('E42').Value = '  -1.85%  '
$ws.Range('E43').Value = '  -3.39%  '
$ws.Range('D44').Value = '15.13'
$ws.Range('E44').Value = '  -6.14%  '
$ws.Range('D45').Value = '0.6572'
$ws.Range('E45').Value = '  -0.98%  '
$ws.Range('D46').Value = '3.999'
$ws.Range('E46').Value = '  +0.82%  '
$ws.Range('D47').Value = '2.272'
$ws.Range('E47').Value = '  -1.85%  '
$ws.Range('D48').Value = '131.91'
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('D49').Value = '0.07924'
$ws.Range('E49').Value = '  -0.79%  '
$ws.Range('D50').Value = '1.218'
$ws.Range('E50').Value = '  +3.93%  '
$ws.Range('D51').Value = '1.231'
$ws.Range('E51').Value = '  +3.56%  '
